$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 2061
$ws.Range("J51").Value = 2061
$ws.Range("L51").Value = 2061
$ws.Range("N51").Value = -3029
# Row 52
$ws.Range("H52").Value = 1249.5
$ws.Range("I52").Value = 1249.5
$ws.Range("K52").Value = 3748.5
$ws.Range("M52").Value = -3588.5
# Row 55
$ws.Range("H55").Value = 143.26315
$ws.Range("I55").Value = 107.36364
$ws.Range("J55").Value = 192.625
$ws.Range("K55").Value = 107.36364
$ws.Range("L55").Value = 192.625
$ws.Range("M55").Value = 106.63636
$ws.Range("N55").Value = -620.625
# Row 103
$ws.Range("H103").Value = 10000636
$ws.Range("I103").Value = 440.44446
$ws.Range("K103").Value = 1321.33338
$ws.Range("M103").Value = -735.33338
# Row 132
$ws.Range("H132").Value = 2734.7715
$ws.Range("I132").Value = 1833.7407
$ws.Range("J132").Value = 5775.75
$ws.Range("K132").Value = 5501.2221
$ws.Range("L132").Value = 17327.25
$ws.Range("M132").Value = -2971.2221
$ws.Range("N132").Value = -22387.25
# Row 138
$ws.Range("H138").Value = 2399.5957
$ws.Range("I138").Value = 2031.0555
$ws.Range("J138").Value = 2628.3447
$ws.Range("K138").Value = 6093.166499999999
$ws.Range("L138").Value = 7885.034100000001
$ws.Range("M138").Value = -953.1664999999994
$ws.Range("N138").Value = -18165.0341

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15018.7
$ws.Range("I32").Value = 4701.0938
$ws.Range("J32").Value = 33361.11
$ws.Range("K32").Value = 4701.0938
$ws.Range("L32").Value = 33361.11
$ws.Range("M32").Value = -4414.0938
$ws.Range("N32").Value = -33935.11
# Row 44
$ws.Range("H44").Value = 34949
$ws.Range("J44").Value = 34949
$ws.Range("L44").Value = 34949
$ws.Range("N44").Value = -35925
# Row 55
$ws.Range("H55").Value = 33653
$ws.Range("J55").Value = 33653
$ws.Range("L55").Value = 33653
$ws.Range("N55").Value = -34283
# Row 76
$ws.Range("H76").Value = 29962.666
$ws.Range("J76").Value = 29962.666
$ws.Range("L76").Value = 29962.666
$ws.Range("N76").Value = -30638.666
# Row 79
$ws.Range("H79").Value = 29962.666
$ws.Range("J79").Value = 29962.666
$ws.Range("L79").Value = 29962.666
$ws.Range("N79").Value = -32302.666
# Row 80
$ws.Range("H80").Value = 42361.6
$ws.Range("J80").Value = 42361.6
$ws.Range("L80").Value = 42361.6
$ws.Range("N80").Value = -44357.6
# Row 83
$ws.Range("H83").Value = 42361.6
$ws.Range("J83").Value = 42361.6
$ws.Range("L83").Value = 127084.8
$ws.Range("N83").Value = -137068.8

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 31758
$ws.Range("J35").Value = 31758
$ws.Range("L35").Value = 31758
$ws.Range("N35").Value = -32378
# Row 51
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 82
$ws.Range("H82").Value = 30513.883
$ws.Range("I82").Value = 11150
$ws.Range("J82").Value = 34663.285
$ws.Range("K82").Value = 11150
$ws.Range("L82").Value = 34663.285
$ws.Range("M82").Value = -10767
$ws.Range("N82").Value = -35429.285
# Row 85
$ws.Range("H85").Value = 30513.883
$ws.Range("I85").Value = 11150
$ws.Range("J85").Value = 34663.285
$ws.Range("K85").Value = 11150
$ws.Range("L85").Value = 34663.285
$ws.Range("M85").Value = -9824
$ws.Range("N85").Value = -37315.285
# Row 134
$ws.Range("H134").Value = 1356.129
$ws.Range("I134").Value = 1360.3448
$ws.Range("J134").Value = 1295
$ws.Range("K134").Value = 4081.0344
$ws.Range("L134").Value = 3885
$ws.Range("M134").Value = -1546.0344
$ws.Range("N134").Value = -8955

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 41
$ws.Range("H41").Value = 20349.75
$ws.Range("J41").Value = 20349.75
$ws.Range("L41").Value = 20349.75
$ws.Range("N41").Value = -21205.75
# Row 51
$ws.Range("H51").Value = 9337.143
$ws.Range("J51").Value = 9337.143
$ws.Range("L51").Value = 9337.143
$ws.Range("N51").Value = -10809.143
# Row 55
$ws.Range("H55").Value = 6191.8
$ws.Range("I55").Value = 6999
$ws.Range("J55").Value = 5990
$ws.Range("K55").Value = 6999
$ws.Range("L55").Value = 5990
$ws.Range("M55").Value = -6684
$ws.Range("N55").Value = -6620
# Row 60
$ws.Range("H60").Value = 30507.111
$ws.Range("J60").Value = 30507.111
$ws.Range("L60").Value = 30507.111
$ws.Range("N60").Value = -31529.111
# Row 61
$ws.Range("H61").Value = 9337.143
$ws.Range("J61").Value = 9337.143
$ws.Range("L61").Value = 9337.143
$ws.Range("N61").Value = -10033.143
# Row 68
$ws.Range("H68").Value = 16907.6
$ws.Range("J68").Value = 16907.6
$ws.Range("L68").Value = 16907.6
$ws.Range("N68").Value = -18405.6
# Row 71
$ws.Range("H71").Value = 16907.6
$ws.Range("J71").Value = 16907.6
$ws.Range("L71").Value = 50722.8
$ws.Range("N71").Value = -58210.8
# Row 132
$ws.Range("H132").Value = 1956.7826
$ws.Range("I132").Value = 2129.5386
$ws.Range("J132").Value = 1732.2
$ws.Range("K132").Value = 6388.6158
$ws.Range("L132").Value = 5196.6
$ws.Range("M132").Value = -3858.6158
$ws.Range("N132").Value = -10256.6

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 632.8570999999999
$ws.Range("I33").Value = 570.25
$ws.Range("K33").Value = 3421.5
$ws.Range("M33").Value = -3138.5
# Row 75
$ws.Range("H75").Value = 6575.625
$ws.Range("J75").Value = 6575.625
$ws.Range("L75").Value = 19726.875
$ws.Range("N75").Value = -21722.875
# Row 78
$ws.Range("H78").Value = 6575.625
$ws.Range("J78").Value = 6575.625
$ws.Range("L78").Value = 59180.625
$ws.Range("N78").Value = -69164.625
# Row 102
$ws.Range("H102").Value = 5922.222
$ws.Range("J102").Value = 5922.222
$ws.Range("L102").Value = 17766.666
$ws.Range("N102").Value = -22634.666

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 3700
$ws.Range("I19").Value = 2200
$ws.Range("J19").Value = 5950
$ws.Range("K19").Value = 2200
$ws.Range("L19").Value = 5950
$ws.Range("M19").Value = -1912
$ws.Range("N19").Value = -6526
# Row 46
$ws.Range("H46").Value = 10988.571
$ws.Range("J46").Value = 12485.454
$ws.Range("L46").Value = 12485.454
$ws.Range("N46").Value = -12797.454
# Row 122
$ws.Range("H122").Value = 2009.7
$ws.Range("I122").Value = 2801.4
$ws.Range("J122").Value = 1218
$ws.Range("K122").Value = 8404.200000000001
$ws.Range("L122").Value = 3654
$ws.Range("M122").Value = -5954.200000000001
$ws.Range("N122").Value = -8554
# Row 123
$ws.Range("H123").Value = 34284.09
$ws.Range("J123").Value = 34284.09
$ws.Range("L123").Value = 34284.09
$ws.Range("N123").Value = -39184.09
# Row 132
$ws.Range("H132").Value = 2071.4358
$ws.Range("I132").Value = 1913.1892
$ws.Range("K132").Value = 5739.5676
$ws.Range("M132").Value = -3209.5676

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1869.6957
$ws.Range("I46").Value = 1726.7333
$ws.Range("J46").Value = 2137.75
$ws.Range("K46").Value = 1726.7333
$ws.Range("L46").Value = 2137.75
$ws.Range("M46").Value = -1538.7333
$ws.Range("N46").Value = -2513.75
# Row 132
$ws.Range("H132").Value = 8328.214
$ws.Range("I132").Value = 9327.091
$ws.Range("J132").Value = 4665.6665
$ws.Range("K132").Value = 27981.273
$ws.Range("L132").Value = 13996.9995
$ws.Range("M132").Value = -25451.273
$ws.Range("N132").Value = -19056.9995

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 19618
$ws.Range("J109").Value = 19618
$ws.Range("L109").Value = 19618
$ws.Range("N109").Value = -22392
